# Updated cryptos list on Fri Jun  9 04:55:31 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Price, Volume(1h)) updates for column D / E.
# Price values are prefixed with a leading single-quote so Excel stores
# them as literal text (matching the workbook's original inline-string
# cells) instead of auto-converting look-alike numbers (e.g. "4.160"
# would otherwise become 4.16, "0.06790" would become 0.0679, etc.)
$updates = @{
    2  = @("26.502.64", "  +0.40%  ")
    3  = @("1.839.56",  "  +0.08%  ")
    4  = @($null,       "  +0.07%  ")
    5  = @("260.18",    "  +0.07%  ")
    7  = @("0.5251",    "  +0.65%  ")
    8  = @("0.3195",    "  -1.17%  ")
    9  = @("0.06790",   "  +0.32%  ")
    10 = @("18.76",     "  +0.78%  ")
    11 = @("0.7847",    "  +2.83%  ")
    12 = @("0.07755",   "  +1.04%  ")
    13 = @("1.835.30",  "  -0.19%  ")
    14 = @("87.87",     "  -0.71%  ")
    15 = @("5.013",     "  -0.05%  ")
    16 = @("1.001",     "  +0.02%  ")
    17 = @("13.84",     "  -0.64%  ")
    19 = @("0.000007947","  +0.53%  ")
    20 = @("26.515.47", "  +0.27%  ")
    21 = @("2.073.16",  "  -0.11%  ")
    22 = @($null,       "  +1.34%  ")
    23 = @("5.967",     "  +0.50%  ")
    24 = @("9.351",     "  -0.96%  ")
    25 = @("141.85",    "  -2.03%  ")
    26 = @("2.187",     "  -1.65%  ")
    27 = @("1.685",     "  +1.68%  ")
    28 = @("16.93",     "  +0.00%  ")
    29 = @("111.74",    "  +0.53%  ")
    30 = @("4.160",     "  -0.46%  ")
    31 = @("0.08699",   "  -0.45%  ")
    32 = @("4.076",     "  -1.48%  ")
    33 = @("0.04882",   "  +1.27%  ")
    34 = @("0.7267",    "  +3.66%  ")
    35 = @("1.136",     "  +1.36%  ")
    36 = @("2.872",     "  +0.99%  ")
    37 = @("3.096",     "  +1.15%  ")
    38 = @("2.244",     "  +3.14%  ")
    39 = @("0.01755",   "  -0.43%  ")
    40 = @("0.4785",    "  -0.83%  ")
    41 = @("0.8944",    "  +0.60%  ")
    42 = @("109.68",    "  -1.47%  ")
    43 = @("5.938",     "  -2.44%  ")
    44 = @($null,       "  +0.14%  ")
    45 = @("7.681",     "  +0.67%  ")
    46 = @("0.4167",    "  +1.33%  ")
    49 = @("0.1233",    "  +1.30%  ")
    50 = @("34.85",     "  +0.49%  ")
    51 = @("0.8912",    "  +1.17%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $price = $vals[0]
    $volume = $vals[1]
    if ($null -ne $price) {
        $ws.Range("D$row").Value = "'" + $price
    }
    $ws.Range("E$row").Value = $volume
}

# Rows 47/48: EnergySwap and Cronos swap rank order (coin name + link),
# each also getting a refreshed price / volume figure.
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05848"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.943"
$ws.Range("E48").Value = "  -0.69%  "
